# Generate Report for Handoff
# Updates the "b.md" rows across the Overview / zh-cn / de-de sheets to
# reflect that b.md is now ready for handoff (a new handoff xlf was
# generated), and records the "not latest" error detail on the language
# sheets.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4bf4307ed0784db95503def832880f2cccd469da/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afffb2700d07e61f945198c75f98c571fcc43852/e2e/b.md."

# ---- Overview sheet: row 3 is the b.md summary row ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusReady
$wsOverview.Range("F3").Value = $statusReady
$wsOverview.Range("G3").Value = "2016-08-22 18:40:13"

# Excel pads ColumnWidth by the default-font "0 digit" margin (~0.8333
# chars for Calibri 11) when it writes the raw OOXML <col width>. Back
# that padding out so the saved width attribute lands on exactly 40.
$columnWidthPad = 0.8333333333333334
$targetColWidth = 40 - $columnWidthPad

# ---- zh-cn sheet: row 3 is the b.md detail row ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusReady
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-22 18:39:57"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = $targetColWidth

# ---- de-de sheet: row 3 is the b.md detail row ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusReady
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-22 18:40:13"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = $targetColWidth
